$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: copy the existing header style (from A1) onto AD1:AF1
# then set the header text values.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record columns (Wins=88, Losses=74, Ties=1) for every player row.
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 1
}
